$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "didn't attend" note to the three students who were marked
# Attendance = "No" / Credit = 0, in column D.
$note = "Didn't attend the lab, didn't show TA the result. "

$ws.Range("D14").Value = $note
$ws.Range("D23").Value = $note
$ws.Range("D30").Value = $note

# Update the active view/selection to match the saved state (scrolled back
# to top, selection on H16 instead of H34).
$ws.Activate()
$ws.Range("H16").Select()
$excel.ActiveWindow.ScrollRow = 1
